$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Post Treatment" (column D) outcome measures
$ws.Range("D2").Value = "Moderately stressful"
$ws.Range("D3").Value = "Moderately stressful"
$ws.Range("D4").Value = "Not stressful"
$ws.Range("D5").Value = "Not stressful"
$ws.Range("D6").Value = "Moderately stressful"
$ws.Range("D7").Value = "A little stressful "

# Widen column D to fit the new content
$ws.Columns.Item(4).ColumnWidth = 20.8333333333333

# Update the active selection to D8, matching the final saved cursor position
$ws.Range("D8").Select()
